$wb = $excel.ActiveWorkbook

# Update both the "展览" sheet and the "全部类型" sheet which both contain
# the same data table: F2 (想去人数) goes 479 -> 480 and F3 goes 62 -> 63.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 480
    $ws.Range("F3").Value = 63
}
